$wb = $excel.ActiveWorkbook

# --- Add the new worksheet "TwoCellsSepartedBySpace" after the last existing sheet ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "TwoCellsSepartedBySpace"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move([System.Reflection.Missing]::Value, $lastSheet)

# Re-fetch a fresh handle by name (Move()/Add() invalidate stale references)
$target = $wb.Worksheets.Item("TwoCellsSepartedBySpace")

# Put "cell1" and "cell2" side by side - two cells separated by a space in the markup
$target.Range("A1").Value = "cell1"
$target.Range("B1").Value = "cell2"

# Make the new sheet the active tab with B1 selected, like the source workbook
$target.Activate()
$target.Range("B1").Select() | Out-Null

# --- Token sheet is no longer the active tab; it also grew an extra (blank) row 6 ---
$tokenSheet = $wb.Worksheets.Item("Token")
$tokenSheet.Range("A2").Copy()
$tokenSheet.Range("A6").PasteSpecial(-4122)  # xlPasteFormats - registers an empty A6 cell
$tokenSheet.Rows.Item(6).RowHeight = 12.35
